$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29 used to be one big (empty) merged cell A29:E29; unmerge it so the
#     new "thrombocytes" entry can carry per-column values/styles like every
#     other row in the table. Columns D/E keep the (already-correct) numeric
#     formatting that the merged cell had, so they are left alone. ---
$ws.Range("A29:E29").UnMerge()

# --- Drop the stray remark in F3 ("actually venous pH for now"); the whole
#     cell goes away (not just its text), same as the other rows that have no
#     "Remarks" entry. ---
$ws.Range("F3").ClearContents()

# --- Pick up the plain-text look used by the other row labels/short names
#     (columns A and B) from an existing row. ---
$ws.Range("A27").Copy()
$ws.Range("A29").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B27").Copy()
$ws.Range("B29").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the new "thrombocytes" variable row. ---
$ws.Range("A29").Value = "thrombocytes"
$ws.Range("B29").Value = "thrombocytes"
$ws.Range("C29").Value = "G/l"
$ws.Range("C29").NumberFormat = "0"
$ws.Range("C29").WrapText = $false
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 1000
$ws.Range("E29").NumberFormat = "#,##0"

# --- Restore the cursor position captured the last time the workbook was
#     saved. ---
$ws.Range("B32").Select()
